$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the latest cryptos data pull.
# D-column values are entered with a leading apostrophe + style reset to keep
# them as plain text (matching the source data, which stores prices as text,
# e.g. "3.062.66") instead of letting Excel auto-convert numeric-looking text
# into a Number cell.
$ws.Range("D2").Value = "'89.423.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.44%  "
$ws.Range("D3").Value = "'3.062.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.98%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'233.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.34%  "
$ws.Range("D6").Value = "'617.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.36%  "
$ws.Range("E7").Value = "  -6.95%  "
$ws.Range("E8").Value = "  -1.89%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "'3.062.78"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.91%  "
$ws.Range("D11").Value = "'0.704"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.42%  "
$ws.Range("E12").Value = "  -1.52%  "
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").Value = "'34.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.93%  "
$ws.Range("D15").Value = "'89.185.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("D16").Value = "'5.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.41%  "
$ws.Range("D17").Value = "'3.637.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.72%  "
$ws.Range("D18").Value = "'3.077.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.02%  "
$ws.Range("D19").Value = "'3.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").Value = "'0.0000213"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").Value = "'13.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.09%  "
$ws.Range("D22").Value = "'428.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.71%  "
$ws.Range("D23").Value = "'5.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.53%  "
$ws.Range("D24").Value = "'8.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.89%  "
$ws.Range("D25").Value = "'5.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.67%  "
$ws.Range("D26").Value = "'86.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -10.40%  "
$ws.Range("E27").Value = "  -6.17%  "
$ws.Range("D28").Value = "'3.237.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.53%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  +11.98%  "
$ws.Range("D31").Value = "'8.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.42%  "
$ws.Range("E32").Value = "  -5.68%  "
$ws.Range("E33").Value = "  -10.81%  "
$ws.Range("D34").Value = "'25.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.20%  "
$ws.Range("E35").Value = "  +2.61%  "
$ws.Range("D36").Value = "'3.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +62.59%  "
$ws.Range("D37").Value = "'6.95"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.13%  "
$ws.Range("D38").Value = "'487.05"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.25%  "
$ws.Range("D39").Value = "'3.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.81%  "
$ws.Range("E40").Value = "  -3.82%  "
$ws.Range("E41").Value = "  -7.55%  "
$ws.Range("E42").Value = "  -2.18%  "
$ws.Range("D43").Value = "'22.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("E45").Value = "  -8.17%  "
$ws.Range("D46").Value = "'157.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.39%  "
$ws.Range("D47").Value = "'1.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.05%  "
$ws.Range("D48").Value = "'0.668"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.80%  "
$ws.Range("D49").Value = "'44.13"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.49%  "
$ws.Range("D50").Value = "'0.998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("E51").Value = "  -6.54%  "
